# Add a "release_date" column (E) to Sheet1, capturing the AACR release
# date associated with each cohort/version block, so that the most recent
# version of a dataset can be determined reliably (NSCLC 2.0 was released
# after NSCLC 2.1, so sorting by version name alone is not reliable).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("E1").Value = "release_date"
$ws.Range("E1").Font.Bold = $true

# Populate the first occurrence of each distinct release date in the same
# order the values were originally entered, so the shared-string table
# ends up in the expected order.
$ws.Range("E13").Value = "2021-08"   # NSCLC v2.1-consortium
$ws.Range("E24").Value = "2022-05"   # NSCLC v2.0-public
$ws.Range("E35").Value = "2021-02"   # CRC v1.1-consortium
$ws.Range("E59").Value = "2021-10"   # BrCa v1.1-consortium
$ws.Range("E2").Value  = "2020-10"   # NSCLC v1.1-consortium

# Fill in the rest of each cohort/version block with its release date.
$ws.Range("E3:E12").Value  = "2020-10"   # NSCLC v1.1-consortium
$ws.Range("E14:E23").Value = "2021-08"   # NSCLC v2.1-consortium
$ws.Range("E25:E34").Value = "2022-05"   # NSCLC v2.0-public
$ws.Range("E36:E46").Value = "2021-02"   # CRC v1.1-consortium
$ws.Range("E47:E58").Value = "2021-08"   # CRC v1.2-consortium
$ws.Range("E60:E70").Value = "2021-10"   # BrCa v1.1-consortium

# Restore the view: scroll back to the top and select G9.
$ws.Activate()
$ws.Range("G9").Select()

Write-Host "release_date column added"
